# Applies the "Update set from data objects" edit: adds 4 new A5InteriorFruits
# breakdown blocks (columns AE:AI, AK:AO, AQ:AU, AW:BA) mirroring the existing
# breakdown blocks in columns G:K / M:Q / S:W, and completes the previously
# partial block in columns Y:AC (filling in AB/AC).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - new group headers (A5InteriorFruits_5..8)
$ws.Range("AE1").Value = 'A5InteriorFruits_5'
$ws.Range("AK1").Value = 'A5InteriorFruits_6'
$ws.Range("AQ1").Value = 'A5InteriorFruits_7'
$ws.Range("AW1").Value = 'A5InteriorFruits_8'

# Row 2 - sub-headers (Overall / Not an area of strength / I am doing okay in this area / This is an area of strength / This is an area of great strength)
$ws.Range("AB2").Value = 'This is an area of strength'
$ws.Range("AC2").Value = 'This is an area of great strength'
$ws.Range("AE2").Value = 'Overall'
$ws.Range("AF2").Value = 'Not an area of strength'
$ws.Range("AG2").Value = 'I am doing okay in this area'
$ws.Range("AH2").Value = 'This is an area of strength'
$ws.Range("AI2").Value = 'This is an area of great strength'
$ws.Range("AK2").Value = 'Overall'
$ws.Range("AL2").Value = 'Not an area of strength'
$ws.Range("AM2").Value = 'I am doing okay in this area'
$ws.Range("AN2").Value = 'This is an area of strength'
$ws.Range("AO2").Value = 'This is an area of great strength'
$ws.Range("AQ2").Value = 'Overall'
$ws.Range("AR2").Value = 'Not an area of strength'
$ws.Range("AS2").Value = 'I am doing okay in this area'
$ws.Range("AT2").Value = 'This is an area of strength'
$ws.Range("AU2").Value = 'This is an area of great strength'
$ws.Range("AW2").Value = 'Overall'
$ws.Range("AX2").Value = 'Not an area of strength'
$ws.Range("AY2").Value = 'I am doing okay in this area'
$ws.Range("AZ2").Value = 'This is an area of strength'
$ws.Range("BA2").Value = 'This is an area of great strength'

# Row 5 - n-size counts
$ws.Range("AB5").Value = 1377
$ws.Range("AC5").Value = 1007
$ws.Range("AE5").Value = 3137
$ws.Range("AF5").Value = 154
$ws.Range("AG5").Value = 974
$ws.Range("AH5").Value = 1349
$ws.Range("AI5").Value = 660
$ws.Range("AK5").Value = 3137
$ws.Range("AL5").Value = 104
$ws.Range("AM5").Value = 1023
$ws.Range("AN5").Value = 1424
$ws.Range("AO5").Value = 586
$ws.Range("AQ5").Value = 3137
$ws.Range("AR5").Value = 359
$ws.Range("AS5").Value = 1606
$ws.Range("AT5").Value = 972
$ws.Range("AU5").Value = 200
$ws.Range("AW5").Value = 3137
$ws.Range("AX5").Value = 512
$ws.Range("AY5").Value = 1728
$ws.Range("AZ5").Value = 744
$ws.Range("BA5").Value = 153

# Row 8 - A2_1: Bible study
$ws.Range("AB8").Value = 0.5453885257806826
$ws.Range("AC8").Value = 0.5163853028798411
$ws.Range("AE8").Value = 0.5259802358941664
$ws.Range("AF8").Value = 0.487012987012987
$ws.Range("AG8").Value = 0.5143737166324436
$ws.Range("AH8").Value = 0.5174203113417346
$ws.Range("AI8").Value = 0.5696969696969697
$ws.Range("AK8").Value = 0.5259802358941664
$ws.Range("AL8").Value = 0.4423076923076923
$ws.Range("AM8").Value = 0.5083088954056696
$ws.Range("AN8").Value = 0.5428370786516854
$ws.Range("AO8").Value = 0.5307167235494881
$ws.Range("AQ8").Value = 0.5259802358941664
$ws.Range("AR8").Value = 0.46518105849582175
$ws.Range("AS8").Value = 0.5236612702366127
$ws.Range("AT8").Value = 0.5473251028806584
$ws.Range("AU8").Value = 0.55
$ws.Range("AW8").Value = 0.5259802358941664
$ws.Range("AX8").Value = 0.529296875
$ws.Range("AY8").Value = 0.5353009259259259
$ws.Range("AZ8").Value = 0.5174731182795699
$ws.Range("BA8").Value = 0.45098039215686275

# Row 9 - A2_2: Small group / ministry / community
$ws.Range("AB9").Value = 0.49963689179375453
$ws.Range("AC9").Value = 0.5114200595829196
$ws.Range("AE9").Value = 0.503984698756774
$ws.Range("AF9").Value = 0.44155844155844154
$ws.Range("AG9").Value = 0.48254620123203285
$ws.Range("AH9").Value = 0.5233506300963677
$ws.Range("AI9").Value = 0.5106060606060606
$ws.Range("AK9").Value = 0.503984698756774
$ws.Range("AL9").Value = 0.47115384615384615
$ws.Range("AM9").Value = 0.4965786901270772
$ws.Range("AN9").Value = 0.5168539325842697
$ws.Range("AO9").Value = 0.49146757679180886
$ws.Range("AQ9").Value = 0.503984698756774
$ws.Range("AR9").Value = 0.467966573816156
$ws.Range("AS9").Value = 0.5049813200498132
$ws.Range("AT9").Value = 0.5246913580246914
$ws.Range("AU9").Value = 0.46
$ws.Range("AW9").Value = 0.503984698756774
$ws.Range("AX9").Value = 0.4609375
$ws.Range("AY9").Value = 0.5266203703703703
$ws.Range("AZ9").Value = 0.4959677419354839
$ws.Range("BA9").Value = 0.43137254901960786

# Row 10 - A2_3: Parent
$ws.Range("AB10").Value = 0.4466230936819172
$ws.Range("AC10").Value = 0.46772591857000995
$ws.Range("AE10").Value = 0.448836467963022
$ws.Range("AF10").Value = 0.461038961038961
$ws.Range("AG10").Value = 0.41786447638603696
$ws.Range("AH10").Value = 0.469236471460341
$ws.Range("AI10").Value = 0.45
$ws.Range("AK10").Value = 0.448836467963022
$ws.Range("AL10").Value = 0.4423076923076923
$ws.Range("AM10").Value = 0.4095796676441838
$ws.Range("AN10").Value = 0.46769662921348315
$ws.Range("AO10").Value = 0.4726962457337884
$ws.Range("AQ10").Value = 0.448836467963022
$ws.Range("AR10").Value = 0.45125348189415043
$ws.Range("AS10").Value = 0.4302615193026152
$ws.Range("AT10").Value = 0.47016460905349794
$ws.Range("AU10").Value = 0.49
$ws.Range("AW10").Value = 0.448836467963022
$ws.Range("AX10").Value = 0.42578125
$ws.Range("AY10").Value = 0.4519675925925926
$ws.Range("AZ10").Value = 0.45564516129032256
$ws.Range("BA10").Value = 0.45751633986928103

# Row 11 - A2_4: Non-parent family member
$ws.Range("AB11").Value = 0.1924473493100944
$ws.Range("AC11").Value = 0.16881827209533268
$ws.Range("AE11").Value = 0.1858463500159388
$ws.Range("AF11").Value = 0.21428571428571427
$ws.Range("AG11").Value = 0.19815195071868583
$ws.Range("AH11").Value = 0.16604892512972572
$ws.Range("AI11").Value = 0.2015151515151515
$ws.Range("AK11").Value = 0.1858463500159388
$ws.Range("AL11").Value = 0.17307692307692307
$ws.Range("AM11").Value = 0.18377321603128055
$ws.Range("AN11").Value = 0.19803370786516855
$ws.Range("AO11").Value = 0.1621160409556314
$ws.Range("AQ11").Value = 0.1858463500159388
$ws.Range("AR11").Value = 0.1977715877437326
$ws.Range("AS11").Value = 0.17496886674968867
$ws.Range("AT11").Value = 0.20164609053497942
$ws.Range("AU11").Value = 0.175
$ws.Range("AW11").Value = 0.1858463500159388
$ws.Range("AX11").Value = 0.1796875
$ws.Range("AY11").Value = 0.1892361111111111
$ws.Range("AZ11").Value = 0.1935483870967742
$ws.Range("BA11").Value = 0.13071895424836602

# Row 12 - A2_5: Friend
$ws.Range("AB12").Value = 0.3769063180827887
$ws.Range("AC12").Value = 0.3743793445878848
$ws.Range("AE12").Value = 0.3736053554351291
$ws.Range("AF12").Value = 0.33766233766233766
$ws.Range("AG12").Value = 0.39117043121149897
$ws.Range("AH12").Value = 0.36619718309859156
$ws.Range("AI12").Value = 0.3712121212121212
$ws.Range("AK12").Value = 0.3736053554351291
$ws.Range("AL12").Value = 0.3557692307692308
$ws.Range("AM12").Value = 0.3782991202346041
$ws.Range("AN12").Value = 0.3714887640449438
$ws.Range("AO12").Value = 0.37372013651877134
$ws.Range("AQ12").Value = 0.3736053554351291
$ws.Range("AR12").Value = 0.3342618384401114
$ws.Range("AS12").Value = 0.3823163138231631
$ws.Range("AT12").Value = 0.39094650205761317
$ws.Range("AU12").Value = 0.29
$ws.Range("AW12").Value = 0.3736053554351291
$ws.Range("AX12").Value = 0.369140625
$ws.Range("AY12").Value = 0.38425925925925924
$ws.Range("AZ12").Value = 0.35618279569892475
$ws.Range("BA12").Value = 0.35294117647058826

# Row 13 - A2_6: Teacher
$ws.Range("AB13").Value = 0.14960058097313
$ws.Range("AC13").Value = 0.1628599801390268
$ws.Range("AE13").Value = 0.15301243226012112
$ws.Range("AF13").Value = 0.11038961038961038
$ws.Range("AG13").Value = 0.14784394250513347
$ws.Range("AH13").Value = 0.15048183839881393
$ws.Range("AI13").Value = 0.17575757575757575
$ws.Range("AK13").Value = 0.15301243226012112
$ws.Range("AL13").Value = 0.15384615384615385
$ws.Range("AM13").Value = 0.14271749755620725
$ws.Range("AN13").Value = 0.15098314606741572
$ws.Range("AO13").Value = 0.1757679180887372
$ws.Range("AQ13").Value = 0.15301243226012112
$ws.Range("AR13").Value = 0.12813370473537605
$ws.Range("AS13").Value = 0.14881693648816938
$ws.Range("AT13").Value = 0.16049382716049382
$ws.Range("AU13").Value = 0.195
$ws.Range("AW13").Value = 0.15301243226012112
$ws.Range("AX13").Value = 0.138671875
$ws.Range("AY13").Value = 0.15046296296296297
$ws.Range("AZ13").Value = 0.16129032258064516
$ws.Range("BA13").Value = 0.1895424836601307

# Row 14 - A2_7: Print and/or digital media (e.g., book, podcast)
$ws.Range("AB14").Value = 0.5882352941176471
$ws.Range("AC14").Value = 0.5888778550148958
$ws.Range("AE14").Value = 0.5890978642014664
$ws.Range("AF14").Value = 0.6233766233766234
$ws.Range("AG14").Value = 0.6149897330595483
$ws.Range("AH14").Value = 0.5774647887323944
$ws.Range("AI14").Value = 0.5666666666666667
$ws.Range("AK14").Value = 0.5890978642014664
$ws.Range("AL14").Value = 0.5673076923076923
$ws.Range("AM14").Value = 0.5806451612903226
$ws.Range("AN14").Value = 0.5990168539325843
$ws.Range("AO14").Value = 0.5836177474402731
$ws.Range("AQ14").Value = 0.5890978642014664
$ws.Range("AR14").Value = 0.6239554317548747
$ws.Range("AS14").Value = 0.5815691158156912
$ws.Range("AT14").Value = 0.5997942386831275
$ws.Range("AU14").Value = 0.535
$ws.Range("AW14").Value = 0.5890978642014664
$ws.Range("AX14").Value = 0.62890625
$ws.Range("AY14").Value = 0.5983796296296297
$ws.Range("AZ14").Value = 0.5524193548387096
$ws.Range("BA14").Value = 0.5294117647058824

# Row 15 - A2_8: Prayer
$ws.Range("AB15").Value = 0.8126361655773421
$ws.Range("AC15").Value = 0.8073485600794439
$ws.Range("AE15").Value = 0.8138348740835193
$ws.Range("AF15").Value = 0.7857142857142857
$ws.Range("AG15").Value = 0.797741273100616
$ws.Range("AH15").Value = 0.816160118606375
$ws.Range("AI15").Value = 0.8393939393939394
$ws.Range("AK15").Value = 0.8138348740835193
$ws.Range("AL15").Value = 0.7692307692307693
$ws.Range("AM15").Value = 0.8103616813294232
$ws.Range("AN15").Value = 0.8230337078651685
$ws.Range("AO15").Value = 0.8054607508532423
$ws.Range("AQ15").Value = 0.8138348740835193
$ws.Range("AR15").Value = 0.7910863509749304
$ws.Range("AS15").Value = 0.811332503113325
$ws.Range("AT15").Value = 0.8261316872427984
$ws.Range("AU15").Value = 0.815
$ws.Range("AW15").Value = 0.8138348740835193
$ws.Range("AX15").Value = 0.826171875
$ws.Range("AY15").Value = 0.8125
$ws.Range("AZ15").Value = 0.8104838709677419
$ws.Range("BA15").Value = 0.803921568627451

# Row 16 - A2_9: Clergy / Religious
$ws.Range("AB16").Value = 0.6289034132171387
$ws.Range("AC16").Value = 0.6434955312810328
$ws.Range("AE16").Value = 0.6362766974816704
$ws.Range("AF16").Value = 0.5584415584415584
$ws.Range("AG16").Value = 0.6303901437371663
$ws.Range("AH16").Value = 0.6493699036323203
$ws.Range("AI16").Value = 0.6363636363636364
$ws.Range("AK16").Value = 0.6362766974816704
$ws.Range("AL16").Value = 0.5865384615384616
$ws.Range("AM16").Value = 0.6011730205278593
$ws.Range("AN16").Value = 0.6601123595505618
$ws.Range("AO16").Value = 0.6484641638225256
$ws.Range("AQ16").Value = 0.6362766974816704
$ws.Range("AR16").Value = 0.6267409470752089
$ws.Range("AS16").Value = 0.638854296388543
$ws.Range("AT16").Value = 0.6481481481481481
$ws.Range("AU16").Value = 0.575
$ws.Range("AW16").Value = 0.6362766974816704
$ws.Range("AX16").Value = 0.634765625
$ws.Range("AY16").Value = 0.6412037037037037
$ws.Range("AZ16").Value = 0.6344086021505376
$ws.Range("BA16").Value = 0.5947712418300654

# Row 17 - A2_10: Event / Encounter
$ws.Range("AB17").Value = 0.37254901960784315
$ws.Range("AC17").Value = 0.38828202581926513
$ws.Range("AE17").Value = 0.3704175964297099
$ws.Range("AF17").Value = 0.4090909090909091
$ws.Range("AG17").Value = 0.3490759753593429
$ws.Range("AH17").Value = 0.374351371386212
$ws.Range("AI17").Value = 0.38484848484848483
$ws.Range("AK17").Value = 0.3704175964297099
$ws.Range("AL17").Value = 0.3942307692307692
$ws.Range("AM17").Value = 0.36950146627565983
$ws.Range("AN17").Value = 0.36235955056179775
$ws.Range("AO17").Value = 0.3873720136518771
$ws.Range("AQ17").Value = 0.3704175964297099
$ws.Range("AR17").Value = 0.3788300835654596
$ws.Range("AS17").Value = 0.35678704856787047
$ws.Range("AT17").Value = 0.38580246913580246
$ws.Range("AU17").Value = 0.39
$ws.Range("AW17").Value = 0.3704175964297099
$ws.Range("AX17").Value = 0.373046875
$ws.Range("AY17").Value = 0.3680555555555556
$ws.Range("AZ17").Value = 0.3844086021505376
$ws.Range("BA17").Value = 0.3202614379084967

# Row 18 - A2_11: Class / Talk
$ws.Range("AB18").Value = 0.29121278140885987
$ws.Range("AC18").Value = 0.28500496524329694
$ws.Range("AE18").Value = 0.2843481032833918
$ws.Range("AF18").Value = 0.2987012987012987
$ws.Range("AG18").Value = 0.2669404517453799
$ws.Range("AH18").Value = 0.29132690882134915
$ws.Range("AI18").Value = 0.2924242424242424
$ws.Range("AK18").Value = 0.2843481032833918
$ws.Range("AL18").Value = 0.23076923076923078
$ws.Range("AM18").Value = 0.27663734115347016
$ws.Range("AN18").Value = 0.28441011235955055
$ws.Range("AO18").Value = 0.30716723549488056
$ws.Range("AQ18").Value = 0.2843481032833918
$ws.Range("AR18").Value = 0.2785515320334262
$ws.Range("AS18").Value = 0.28393524283935245
$ws.Range("AT18").Value = 0.28703703703703703
$ws.Range("AU18").Value = 0.285
$ws.Range("AW18").Value = 0.2843481032833918
$ws.Range("AX18").Value = 0.263671875
$ws.Range("AY18").Value = 0.30324074074074076
$ws.Range("AZ18").Value = 0.25806451612903225
$ws.Range("BA18").Value = 0.2679738562091503

# Row 19 - A2_12: Sacrament
$ws.Range("AB19").Value = 0.7458242556281772
$ws.Range("AC19").Value = 0.7576961271102284
$ws.Range("AE19").Value = 0.7593241950908511
$ws.Range("AF19").Value = 0.7337662337662337
$ws.Range("AG19").Value = 0.7351129363449692
$ws.Range("AH19").Value = 0.7679762787249814
$ws.Range("AI19").Value = 0.7833333333333333
$ws.Range("AK19").Value = 0.7593241950908511
$ws.Range("AL19").Value = 0.7307692307692307
$ws.Range("AM19").Value = 0.729227761485826
$ws.Range("AN19").Value = 0.7710674157303371
$ws.Range("AO19").Value = 0.78839590443686
$ws.Range("AQ19").Value = 0.7593241950908511
$ws.Range("AR19").Value = 0.7325905292479109
$ws.Range("AS19").Value = 0.7615193026151931
$ws.Range("AT19").Value = 0.7664609053497943
$ws.Range("AU19").Value = 0.755
$ws.Range("AW19").Value = 0.7593241950908511
$ws.Range("AX19").Value = 0.763671875
$ws.Range("AY19").Value = 0.7604166666666666
$ws.Range("AZ19").Value = 0.7486559139784946
$ws.Range("BA19").Value = 0.7843137254901961

# Row 20 - A2_13: Lives of the saints
$ws.Range("AB20").Value = 0.439360929557008
$ws.Range("AC20").Value = 0.44985104270109233
$ws.Range("AE20").Value = 0.4430985017532674
$ws.Range("AF20").Value = 0.37012987012987014
$ws.Range("AG20").Value = 0.4271047227926078
$ws.Range("AH20").Value = 0.4425500370644922
$ws.Range("AI20").Value = 0.48484848484848486
$ws.Range("AK20").Value = 0.4430985017532674
$ws.Range("AL20").Value = 0.36538461538461536
$ws.Range("AM20").Value = 0.41055718475073316
$ws.Range("AN20").Value = 0.45997191011235955
$ws.Range("AO20").Value = 0.4726962457337884
$ws.Range("AQ20").Value = 0.4430985017532674
$ws.Range("AR20").Value = 0.42896935933147634
$ws.Range("AS20").Value = 0.43711083437110837
$ws.Range("AT20").Value = 0.44753086419753085
$ws.Range("AU20").Value = 0.495
$ws.Range("AW20").Value = 0.4430985017532674
$ws.Range("AX20").Value = 0.453125
$ws.Range("AY20").Value = 0.44560185185185186
$ws.Range("AZ20").Value = 0.4260752688172043
$ws.Range("BA20").Value = 0.46405228758169936

# Row 21 - A2_14: Other (specify)
$ws.Range("AB21").Value = 0.13870733478576616
$ws.Range("AC21").Value = 0.1529294935451837
$ws.Range("AE21").Value = 0.14026139623844439
$ws.Range("AF21").Value = 0.12337662337662338
$ws.Range("AG21").Value = 0.12833675564681724
$ws.Range("AH21").Value = 0.13936249073387694
$ws.Range("AI21").Value = 0.16363636363636364
$ws.Range("AK21").Value = 0.14026139623844439
$ws.Range("AL21").Value = 0.15384615384615385
$ws.Range("AM21").Value = 0.13978494623655913
$ws.Range("AN21").Value = 0.1306179775280899
$ws.Range("AO21").Value = 0.1621160409556314
$ws.Range("AQ21").Value = 0.14026139623844439
$ws.Range("AR21").Value = 0.13649025069637882
$ws.Range("AS21").Value = 0.13138231631382316
$ws.Range("AT21").Value = 0.15020576131687244
$ws.Range("AU21").Value = 0.17
$ws.Range("AW21").Value = 0.14026139623844439
$ws.Range("AX21").Value = 0.1640625
$ws.Range("AY21").Value = 0.13541666666666666
$ws.Range("AZ21").Value = 0.1303763440860215
$ws.Range("BA21").Value = 0.16339869281045752

# Row 22 - A2_15: Nothing has had a meaningful impact on my faith life
$ws.Range("AE22").Value = 0.000318775900541919
$ws.Range("AG22").Value = 0.001026694045174538
$ws.Range("AK22").Value = 0.000318775900541919
$ws.Range("AN22").Value = 0.0007022471910112359
$ws.Range("AQ22").Value = 0.000318775900541919
$ws.Range("AS22").Value = 0.0006226650062266501
$ws.Range("AW22").Value = 0.000318775900541919
$ws.Range("AZ22").Value = 0.0013440860215053765
